$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3
$ws.Cells.Item(2, 3).Value = 'flower/flower027.jpg'
$ws.Cells.Item(2, 4).Value = 'dienen'
$ws.Cells.Item(2, 5).Value = 'flower'

$ws.Cells.Item(3, 2).Value = 70
$ws.Cells.Item(3, 3).Value = 'face/face031.jpg'
$ws.Cells.Item(3, 4).Value = 'passen'
$ws.Cells.Item(3, 5).Value = 'face'

$ws.Cells.Item(4, 2).Value = 41
$ws.Cells.Item(4, 3).Value = 'flower/flower024.jpg'
$ws.Cells.Item(4, 4).Value = 'trotzen'
$ws.Cells.Item(4, 5).Value = 'flower'

$ws.Cells.Item(5, 2).Value = 118
$ws.Cells.Item(5, 3).Value = 'face/face013.jpg'
$ws.Cells.Item(5, 4).Value = 'proben'
$ws.Cells.Item(5, 5).Value = 'face'

$ws.Cells.Item(6, 2).Value = 33
$ws.Cells.Item(6, 3).Value = 'face/face015.jpg'
$ws.Cells.Item(6, 4).Value = 'frischen'
$ws.Cells.Item(6, 5).Value = 'face'

$ws.Cells.Item(7, 2).Value = 105
$ws.Cells.Item(7, 3).Value = 'face/face011.jpg'
$ws.Cells.Item(7, 4).Value = 'rufen'
$ws.Cells.Item(7, 5).Value = 'face'

$ws.Cells.Item(8, 2).Value = 126
$ws.Cells.Item(8, 3).Value = 'face/face000.jpg'
$ws.Cells.Item(8, 4).Value = 'leeren'
$ws.Cells.Item(8, 5).Value = 'face'

$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = 'flower/flower030.jpg'
$ws.Cells.Item(9, 4).Value = 'planen'
$ws.Cells.Item(9, 5).Value = 'flower'

$ws.Cells.Item(10, 2).Value = 24
$ws.Cells.Item(10, 3).Value = 'flower/flower018.jpg'
$ws.Cells.Item(10, 4).Value = 'wählen'
$ws.Cells.Item(10, 5).Value = 'flower'

$ws.Cells.Item(11, 2).Value = 92
$ws.Cells.Item(11, 3).Value = 'flower/flower000.jpg'
$ws.Cells.Item(11, 4).Value = 'hacken'
$ws.Cells.Item(11, 5).Value = 'flower'

$ws.Cells.Item(12, 2).Value = 42
$ws.Cells.Item(12, 3).Value = 'flower/flower031.jpg'
$ws.Cells.Item(12, 4).Value = 'spenden'
$ws.Cells.Item(12, 5).Value = 'flower'

$ws.Cells.Item(13, 2).Value = 80
$ws.Cells.Item(13, 3).Value = 'flower/flower016.jpg'
$ws.Cells.Item(13, 4).Value = 'heben'
$ws.Cells.Item(13, 5).Value = 'flower'

$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 'face/face027.jpg'
$ws.Cells.Item(14, 4).Value = 'grenzen'
$ws.Cells.Item(14, 5).Value = 'face'

$ws.Cells.Item(15, 2).Value = 102
$ws.Cells.Item(15, 3).Value = 'flower/flower004.jpg'
$ws.Cells.Item(15, 4).Value = 'lügen'
$ws.Cells.Item(15, 5).Value = 'flower'

$ws.Cells.Item(16, 2).Value = 75
$ws.Cells.Item(16, 3).Value = 'face/face006.jpg'
$ws.Cells.Item(16, 4).Value = 'parken'
$ws.Cells.Item(16, 5).Value = 'face'

$ws.Cells.Item(17, 2).Value = 52
$ws.Cells.Item(17, 3).Value = 'flower/flower023.jpg'
$ws.Cells.Item(17, 4).Value = 'streifen'
$ws.Cells.Item(17, 5).Value = 'flower'

$ws.Cells.Item(18, 2).Value = 85
$ws.Cells.Item(18, 3).Value = 'face/face005.jpg'
$ws.Cells.Item(18, 4).Value = 'angeln'
$ws.Cells.Item(18, 5).Value = 'face'

$ws.Cells.Item(19, 2).Value = 113
$ws.Cells.Item(19, 3).Value = 'face/face007.jpg'
$ws.Cells.Item(19, 4).Value = 'fließen'
$ws.Cells.Item(19, 5).Value = 'face'

$ws.Cells.Item(20, 2).Value = 7
$ws.Cells.Item(20, 3).Value = 'flower/flower021.jpg'
$ws.Cells.Item(20, 4).Value = 'zielen'
$ws.Cells.Item(20, 5).Value = 'flower'

$ws.Cells.Item(21, 2).Value = 12
$ws.Cells.Item(21, 3).Value = 'flower/flower009.jpg'
$ws.Cells.Item(21, 4).Value = 'bergen'
$ws.Cells.Item(21, 5).Value = 'flower'

$ws.Cells.Item(22, 2).Value = 127
$ws.Cells.Item(22, 3).Value = 'face/face001.jpg'
$ws.Cells.Item(22, 4).Value = 'hören'
$ws.Cells.Item(22, 5).Value = 'face'

$ws.Cells.Item(23, 2).Value = 31
$ws.Cells.Item(23, 3).Value = 'face/face029.jpg'
$ws.Cells.Item(23, 4).Value = 'stoßen'
$ws.Cells.Item(23, 5).Value = 'face'

$ws.Cells.Item(24, 2).Value = 29
$ws.Cells.Item(24, 3).Value = 'face/face030.jpg'
$ws.Cells.Item(24, 4).Value = 'nullen'
$ws.Cells.Item(24, 5).Value = 'face'

$ws.Cells.Item(25, 2).Value = 9
$ws.Cells.Item(25, 3).Value = 'face/face002.jpg'
$ws.Cells.Item(25, 4).Value = 'kranken'
$ws.Cells.Item(25, 5).Value = 'face'

$ws.Cells.Item(26, 2).Value = 65
$ws.Cells.Item(26, 3).Value = 'face/face022.jpg'
$ws.Cells.Item(26, 4).Value = 'lernen'
$ws.Cells.Item(26, 5).Value = 'face'

$ws.Cells.Item(27, 2).Value = 97
$ws.Cells.Item(27, 3).Value = 'flower/flower008.jpg'
$ws.Cells.Item(27, 4).Value = 'kriegen'
$ws.Cells.Item(27, 5).Value = 'flower'

$ws.Cells.Item(28, 2).Value = 117
$ws.Cells.Item(28, 3).Value = 'flower/flower002.jpg'
$ws.Cells.Item(28, 4).Value = 'münzen'
$ws.Cells.Item(28, 5).Value = 'flower'

$ws.Cells.Item(29, 2).Value = 19
$ws.Cells.Item(29, 3).Value = 'flower/flower029.jpg'
$ws.Cells.Item(29, 4).Value = 'ändern'
$ws.Cells.Item(29, 5).Value = 'flower'

$ws.Cells.Item(30, 2).Value = 50
$ws.Cells.Item(30, 3).Value = 'face/face021.jpg'
$ws.Cells.Item(30, 4).Value = 'wehen'
$ws.Cells.Item(30, 5).Value = 'face'

$ws.Cells.Item(31, 2).Value = 79
$ws.Cells.Item(31, 3).Value = 'flower/flower026.jpg'
$ws.Cells.Item(31, 4).Value = 'spüren'
$ws.Cells.Item(31, 5).Value = 'flower'

$ws.Cells.Item(32, 2).Value = 16
$ws.Cells.Item(32, 3).Value = 'flower/flower012.jpg'
$ws.Cells.Item(32, 4).Value = 'bauen'
$ws.Cells.Item(32, 5).Value = 'flower'

$ws.Cells.Item(33, 2).Value = 89
$ws.Cells.Item(33, 3).Value = 'face/face008.jpg'
$ws.Cells.Item(33, 4).Value = 'kennen'
$ws.Cells.Item(33, 5).Value = 'face'
